$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This commit regenerates the localization-status report: the old source
# file "26289a28-...md" was re-run and produced a new handoff (new source
# guid "90c865fa-...md", new transform hash, new handoff timestamps), and
# the previously-failed-transform row ("3ff4ef4e-...md") is gone (dropped)
# from every sheet - the ".localization-config" row moves up to take its
# place as row 3.
# ---------------------------------------------------------------------------

$oldSrc  = "26289a28-d5df-4959-b380-108546cc004e.md"
$newSrc  = "90c865fa-efe0-478f-826f-f969f6876da8.md"

$oldHashZh = "26289a28-d5df-4959-b380-108546cc004e.11d843c29faa8586a7ca71135c242699e814550b.zh-cn.xlf"
$newHashZh = "90c865fa-efe0-478f-826f-f969f6876da8.a37e708a1df312362795fd615e4c978c2ca94acf.zh-cn.xlf"

$oldHashDe = "26289a28-d5df-4959-b380-108546cc004e.11d843c29faa8586a7ca71135c242699e814550b.de-de.xlf"
$newHashDe = "90c865fa-efe0-478f-826f-f969f6876da8.a37e708a1df312362795fd615e4c978c2ca94acf.de-de.xlf"

$newDateZh = "2016-01-13 12:06:05"
$newDateDe = "2016-01-13 12:06:28"

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/185f535b5b613f76de96b6d2442c2e0c38b3e202"

# =============================== Overview ==================================
$ws = $wb.Worksheets.Item("Overview")

# Drop the "Handoff transform failed" row (row 3); the ".localization-config"
# row (old row 4) shifts up into row 3 automatically.
$ws.Rows.Item(3).Delete()

# Row 2 now points at the regenerated source file.
$ws.Range("A2").Value = $newSrc

# The hyperlink table is stale after the row shift (old refs/targets) -
# clear it out and rebuild it against the final two data rows.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newSrc", "", "", $newSrc) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/.localization-config", "", "", ".localization-config") | Out-Null

# =============================== zh-cn ======================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value = $newSrc
$ws.Range("C2").Value = $newHashZh
$ws.Range("D2").Value = $newDateZh

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newSrc", "", "", $newSrc) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4265f07f08632b9f20a3598f7466642cd68837a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$newHashZh", "", "", $newHashZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/.localization-config", "", "", ".localization-config") | Out-Null

# =============================== de-de ======================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value = $newSrc
$ws.Range("C2").Value = $newHashDe
$ws.Range("D2").Value = $newDateDe

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newSrc", "", "", $newSrc) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf07010c034ceb61211d2906c0355f1df1402851/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$newHashDe", "", "", $newHashDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "Report regenerated for handoff"
